# Append two new daily rows (2025-11-27 and 2025-11-28) to each of the
# five data worksheets, carrying forward the existing date/number
# formatting from the preceding row.

$wb = $excel.ActiveWorkbook

# Sheet name -> [ B value for 45988 (11/27), B value for 45989 (11/28) ]
$newData = @{
    "카카오"     = @(721900, 0)
    "NAVER"      = @(1332180, 0)
    "농심"       = @(115414, 0)
    "삼양식품"   = @(439555, 0)
    "엔씨소프트" = @(248197, 0)
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $newData.ContainsKey($name)) { continue }

    $values = $newData[$name]

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $lastDateCell = $ws.Cells.Item($lastRow, 1)

    $row1 = $lastRow + 1
    $row2 = $lastRow + 2

    $aCell1 = $ws.Cells.Item($row1, 1)
    $bCell1 = $ws.Cells.Item($row1, 2)
    $aCell2 = $ws.Cells.Item($row2, 1)
    $bCell2 = $ws.Cells.Item($row2, 2)

    $aCell1.Value = 45988
    $aCell1.NumberFormat = $lastDateCell.NumberFormat
    $bCell1.Value = $values[0]

    $aCell2.Value = 45989
    $aCell2.NumberFormat = $lastDateCell.NumberFormat
    $bCell2.Value = $values[1]
}
